$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. D13: "degree.." -> "degree = default (3)"
$ws.Range("D13").Value = "degree = default (3)"

# 2. Row 6 (SVC / RBF / tweak gamma = 10) data re-entered with new values
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1
$ws.Range("H6").Value = 1
$ws.Range("I6").Value = 1
$ws.Range("K6").Value = 0.50241546000000004
$ws.Range("L6").Value = 0.57971013999999998

# K6:L6 pick up the regular bordered style used elsewhere in the table
# (they previously used the un-bordered "gamma=auto" style)
$ws.Range("K4").Copy()
$ws.Range("K6:L6").PasteSpecial(-4122)

# 3. Bold the "avg training acc" (J) and "avg testing acc" (P) columns
$ws.Range("J4:J15").Font.Bold = $true
$ws.Range("P4:P15").Font.Bold = $true

# 4. New summary row: minimum of the avg-training-acc column
$ws.Range("J16").Copy()
$ws.Range("J17").PasteSpecial(-4122)
$ws.Range("J17").Formula = "=MIN(J4:J15)"
$ws.Range("J17").Font.Bold = $true

# 5. New discussion paragraph, merged & word-wrapped under "Discuss your results"
$discussion = $ws.Range("B19:H23")
$discussion.Merge()
$discussion.WrapText = $true
$discussion.HorizontalAlignment = -4131
$discussion.VerticalAlignment = -4160
$ws.Range("B19").Value = "As we can see In the above table, the highest overall accuracy we got was for SVC 'Linear Kernel' and the worst accuracy we got was from the 'Sigmoid Kernel'. For RBF, the best values if gamma we found was the default one which used the 1 / (n_features * X.var()) as value of gamma. We found out that lower the value of gamma better the accuracy was. For c, the best accuracy was for 10 after which the accuracy didnt improve by much and started to lower again"

# 6. Leave the cursor where the author left it
$ws.Range("D13").Select()
